$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.184.59"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").Value = "2.617.28"
$ws.Range("E3").Value = "  -1.07%  "
$ws.Range("E4").Value = "  +0.01%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "588.95"
$c.ClearFormats()
$ws.Range("E5").Value = "  -1.76%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "165.65"
$c.ClearFormats()
$ws.Range("E6").Value = "  -1.70%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -2.52%  "
$ws.Range("D9").Value = "2.616.88"
$ws.Range("E9").Value = "  -1.08%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.138"
$c.ClearFormats()
$ws.Range("E10").Value = "  -5.05%  "
$ws.Range("E11").Value = "  +0.94%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.364"
$c.ClearFormats()
$ws.Range("E12").Value = "  -0.60%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "5.21"
$c.ClearFormats()
$ws.Range("E13").Value = "  -0.59%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "27.31"
$c.ClearFormats()
$ws.Range("E14").Value = "  -2.69%  "
$ws.Range("E15").Value = "  -0.95%  "
$ws.Range("E16").Value = "  -2.87%  "
$ws.Range("D17").Value = "66.987.90"
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").Value = "2.618.84"
$ws.Range("E18").Value = "  -0.88%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "11.75"
$c.ClearFormats()
$ws.Range("E19").Value = "  -1.38%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "7.79"
$c.ClearFormats()
$ws.Range("E20").Value = "  -1.24%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "355.10"
$c.ClearFormats()
$ws.Range("E21").Value = "  -2.39%  "
$ws.Range("E22").Value = "  -3.11%  "
$ws.Range("E23").Value = "  -3.26%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "10.50"
$c.ClearFormats()
$ws.Range("E24").Value = "  -3.76%  "
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("E26").Value = "  -4.93%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "69.14"
$c.ClearFormats()
$ws.Range("E27").Value = "  -2.41%  "
$ws.Range("D28").Value = "2.751.60"
$ws.Range("E28").Value = "  -1.00%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.ClearFormats()
$ws.Range("E29").Value = "  +0.32%  "
$ws.Range("D30").Value = "0.0₃0996"
$ws.Range("E30").Value = "  -3.06%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "544.81"
$c.ClearFormats()
$ws.Range("E31").Value = "  -2.44%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "7.86"
$c.ClearFormats()
$ws.Range("E32").Value = "  -2.40%  "
$ws.Range("E33").Value = "  -4.31%  "
$ws.Range("E34").Value = "  -2.83%  "
$ws.Range("E35").Value = "  -0.02%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.ClearFormats()
$ws.Range("E37").Value = "  -4.31%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "157.63"
$c.ClearFormats()
$ws.Range("E38").Value = "  -0.28%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "18.90"
$c.ClearFormats()
$ws.Range("E39").Value = "  -2.77%  "
$ws.Range("E40").Value = "  -2.62%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "18.24"
$c.ClearFormats()
$ws.Range("E41").Value = "  +1.70%  "
$ws.Range("E42").Value = "  -1.91%  "
$ws.Range("E43").Value = "  -2.92%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("E45").Value = "  -5.00%  "
$ws.Range("D46").Value = "0.0₆0298"
$ws.Range("E46").Value = "  -1.17%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "151.35"
$c.ClearFormats()
$ws.Range("E47").Value = "  -1.62%  "
$ws.Range("E48").Value = "  -3.56%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "3.77"
$c.ClearFormats()
$ws.Range("E49").Value = "  -3.19%  "
$ws.Range("E50").Value = "  -1.54%  "
$ws.Range("E51").Value = "  -1.30%  "
